# Update cryptos list D (Price) and E (Volume(1h)) columns with latest scrape values.
# D-column values are forced to text (matching the source data, which stores prices
# as literal strings such as "42.812.75") by briefly switching NumberFormat to "@"
# then restoring the cell style to "Normal" so no stray number format lingers.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "42.812.75"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +3.50%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.254.48"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +3.02%  "
$ws.Range("E4").Value = "  -0.10%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "253.82"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.53%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.637"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +1.31%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "72.16"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +5.38%  "
$ws.Range("E8").Value = "  -0.18%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.650"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +13.23%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "40.91"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +8.42%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "59.54"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +0.91%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.0955"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +2.10%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "7.43"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +4.11%  "
$ws.Range("E14").Value = "  -0.47%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "2.590.43"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +2.88%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.888"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +1.96%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "14.72"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +1.41%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "2.264.44"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +2.31%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "42.790.43"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +3.53%  "
$ws.Range("E20").Value = "  +2.07%  "
$ws.Range("E21").Value = "  +1.08%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "73.11"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +1.44%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "236.32"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +1.40%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.12"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +5.05%  "
$ws.Range("E25").Value = "  +0.18%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "11.74"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -0.36%  "
$ws.Range("E27").Value = "  -0.01%  "
$ws.Range("E28").Value = "  -2.87%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "3.71"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -0.26%  "
$ws.Range("E30").Value = "  +2.29%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "167.71"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -0.54%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "21.05"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +1.88%  "
$ws.Range("E33").Value = "  +8.61%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "6.18"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +12.34%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.0786"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +4.04%  "
$ws.Range("E36").Value = "  +1.13%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "29.17"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +8.66%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "4.71"
$ws.Range("D38").Style = "Normal"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "4.13"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -0.61%  "
$ws.Range("E40").Value = "  +8.31%  "
$ws.Range("E41").Value = "  +3.60%  "
$ws.Range("E42").Value = "  +4.48%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "12.55"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +0.84%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "64.31"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -0.03%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "5.01"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -1.67%  "
$ws.Range("E46").Value = "  +4.59%  "
$ws.Range("E47").Value = "  +0.59%  "
$ws.Range("E48").Value = "  +0.69%  "
$ws.Range("E49").Value = "  -1.31%  "
$ws.Range("E50").Value = "  -0.21%  "
$ws.Range("E51").Value = "  +1.43%  "
